$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.301.08'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.21%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.872.43'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.11%  '

$ws.Range('E4').Value = '  +0.03%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7082'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.60%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '241.85'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.11%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.000'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.02%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07793'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.35%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3106'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.11%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '25.08'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.38%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08381'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.11%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.872.21'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.51%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.230'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.17%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.7178'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.54%  '

$ws.Range('E15').Value = '  -0.46%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000008367'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.23%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.134'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.86%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '29.316.04'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.22%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '239.94'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.44%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.127.97'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.06%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.19'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.03%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9998'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.01%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.741'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.93%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.001'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.04%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1589'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.50%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '162.82'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.70%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.022'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.24%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.47'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.51%  '

$ws.Range('E29').Value = '  -0.13%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.425'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.49%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.347'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.48%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.239'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.86%  '

$ws.Range('E33').Value = '  +2.16%  '

$ws.Range('E34').Value = '  +0.97%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7535'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.21%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.174'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.12%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.686'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.21%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01874'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.69%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.238.91'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +6.62%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.731'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.30%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.536'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.05%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8931'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.54%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '109.56'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.77%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '72.26'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.43%  '

$ws.Range('E45').Value = '  +0.01%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000128'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +6.02%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.027.30'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.31%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5200'
$ws.Range('D48').Style = 'Normal'

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.792'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.10%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.443'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.57%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4339'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.91%  '
